$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 11 (pushing the "301..." block down to rows 13+).
# Copy row 10's formatting/values down twice so the new rows inherit the same
# cell styles (fill/font/alignment) as the rest of the "General" (82) block,
# then restore the thin borders that a row-insert drops, and finally
# overwrite with the real content for the two new log-message rows.
$ws.Rows("10:10").Copy()
$ws.Rows("11:11").Insert()
$ws.Rows("10:10").Copy()
$ws.Rows("11:11").Insert()
$ws.Range("A11:F12").Borders.LineStyle = 1

# Row 11: new "warn" entry (code 209)
$ws.Range("A11").Value = "Both console and output log have disabled messages."
$ws.Range("B11").Value = "General"
$ws.Range("C11").Value = 209
$ws.Range("D11").Value = "warn"
$ws.Range("E11").Value = "ElevatorOptions.__init__()"
$ws.Range("F11").Value = ""

# Row 12: new "error" entry (code 210)
$ws.Range("A12").Value = "OSError [message]"
$ws.Range("B12").Value = "General"
$ws.Range("C12").Value = 210
$ws.Range("D12").Value = "error"
$ws.Range("E12").Value = "elevate_file, elevate_string, elevate_package"
$ws.Range("F12").Value = ""

$ws.Range("F12").Select()
